$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "bcsstk20.mtx"
$ws.Cells.Item(2, 3).Value = "MP"
$ws.Cells.Item(2, 4).Value = 12614913497658990
$ws.Cells.Item(2, 5).Value = 30
$ws.Cells.Item(2, 6).Value = 0.006764411926269531
$ws.Cells.Item(2, 7).Value = 485

$ws.Cells.Item(3, 2).Value = "bcsstk20.mtx"
$ws.Cells.Item(3, 3).Value = "MP_Aitken"
$ws.Cells.Item(3, 4).Value = 12615279568909070
$ws.Cells.Item(3, 5).Value = 20
$ws.Cells.Item(3, 6).Value = 0.001083612442016602
$ws.Cells.Item(3, 7).Value = 485

$ws.Cells.Item(4, 2).Value = "bcsstk25.mtx"
$ws.Cells.Item(4, 3).Value = "MP"
$ws.Cells.Item(4, 4).Value = 1060008038118777
$ws.Cells.Item(4, 5).Value = 23
$ws.Cells.Item(4, 6).Value = 2.791603803634644
$ws.Cells.Item(4, 7).Value = 15439

$ws.Cells.Item(5, 2).Value = "bcsstk25.mtx"
$ws.Cells.Item(5, 3).Value = "MP_Aitken"
$ws.Cells.Item(5, 4).Value = 1060022565795500
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 1.924396276473999
$ws.Cells.Item(5, 7).Value = 15439

$ws.Cells.Item(6, 2).Value = "bcsstk22.mtx"
$ws.Cells.Item(6, 3).Value = "MP"
$ws.Cells.Item(6, 4).Value = 5849890.36991932
$ws.Cells.Item(6, 5).Value = 32
$ws.Cells.Item(6, 6).Value = 0.002159357070922852
$ws.Cells.Item(6, 7).Value = 138

$ws.Cells.Item(7, 2).Value = "bcsstk22.mtx"
$ws.Cells.Item(7, 3).Value = "MP_Aitken"
$ws.Cells.Item(7, 4).Value = 5849906.436342365
$ws.Cells.Item(7, 5).Value = 17
$ws.Cells.Item(7, 6).Value = 0.001317739486694336
$ws.Cells.Item(7, 7).Value = 138

$ws.Cells.Item(8, 2).Value = "bcsstk26.mtx"
$ws.Cells.Item(8, 3).Value = "MP"
$ws.Cells.Item(8, 4).Value = 158268836382.195
$ws.Cells.Item(8, 5).Value = 46
$ws.Cells.Item(8, 6).Value = 0.08556103706359863
$ws.Cells.Item(8, 7).Value = 1922

$ws.Cells.Item(9, 2).Value = "bcsstk26.mtx"
$ws.Cells.Item(9, 3).Value = "MP_Aitken"
$ws.Cells.Item(9, 4).Value = 158274519528.784
$ws.Cells.Item(9, 5).Value = 34
$ws.Cells.Item(9, 6).Value = 0.06520462036132812
$ws.Cells.Item(9, 7).Value = 1922

$ws.Cells.Item(10, 2).Value = "bcsstk27.mtx"
$ws.Cells.Item(10, 3).Value = "MP"
$ws.Cells.Item(10, 4).Value = 3464617.500374258
$ws.Cells.Item(10, 5).Value = 85
$ws.Cells.Item(10, 6).Value = 0.05857491493225098
$ws.Cells.Item(10, 7).Value = 1224

$ws.Cells.Item(11, 2).Value = "bcsstk27.mtx"
$ws.Cells.Item(11, 3).Value = "MP_Aitken"
$ws.Cells.Item(11, 4).Value = 3464838.462803914
$ws.Cells.Item(11, 5).Value = 49
$ws.Cells.Item(11, 6).Value = 0.02513456344604492
$ws.Cells.Item(11, 7).Value = 1224

$ws.Cells.Item(12, 2).Value = "bcsstk21.mtx"
$ws.Cells.Item(12, 3).Value = "MP"
$ws.Cells.Item(12, 4).Value = 127119840.5799289
$ws.Cells.Item(12, 5).Value = 335
$ws.Cells.Item(12, 6).Value = 2.335164546966553
$ws.Cells.Item(12, 7).Value = 3600

$ws.Cells.Item(13, 2).Value = "bcsstk21.mtx"
$ws.Cells.Item(13, 3).Value = "MP_Aitken"
$ws.Cells.Item(13, 4).Value = 127191447.3212106
$ws.Cells.Item(13, 5).Value = 185
$ws.Cells.Item(13, 6).Value = 1.308945417404175
$ws.Cells.Item(13, 7).Value = 3600

$ws.Cells.Item(14, 2).Value = "bcsstk19.mtx"
$ws.Cells.Item(14, 3).Value = "MP"
$ws.Cells.Item(14, 4).Value = 192111739442336.5
$ws.Cells.Item(14, 5).Value = 23
$ws.Cells.Item(14, 6).Value = 0.005386829376220703
$ws.Cells.Item(14, 7).Value = 817

$ws.Cells.Item(15, 2).Value = "bcsstk19.mtx"
$ws.Cells.Item(15, 3).Value = "MP_Aitken"
$ws.Cells.Item(15, 4).Value = 192109842633671.2
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 0.003606319427490234
$ws.Cells.Item(15, 7).Value = 817

$ws.Cells.Item(16, 2).Value = "bcsstk23.mtx"
$ws.Cells.Item(16, 3).Value = "MP"
$ws.Cells.Item(16, 4).Value = 22579879913205152
$ws.Cells.Item(16, 5).Value = 31
$ws.Cells.Item(16, 6).Value = 0.1412050724029541
$ws.Cells.Item(16, 7).Value = 3134

$ws.Cells.Item(17, 2).Value = "bcsstk23.mtx"
$ws.Cells.Item(17, 3).Value = "MP_Aitken"
$ws.Cells.Item(17, 4).Value = 22580409268915860
$ws.Cells.Item(17, 5).Value = 19
$ws.Cells.Item(17, 6).Value = 0.09653830528259277
$ws.Cells.Item(17, 7).Value = 3134

$ws.Cells.Item(18, 2).Value = "bcsstk17.mtx"
$ws.Cells.Item(18, 3).Value = "MP"
$ws.Cells.Item(18, 4).Value = 12960385068.16294
$ws.Cells.Item(18, 5).Value = 22
$ws.Cells.Item(18, 6).Value = 1.335657835006714
$ws.Cells.Item(18, 7).Value = 10974

$ws.Cells.Item(19, 2).Value = "bcsstk17.mtx"
$ws.Cells.Item(19, 3).Value = "MP_Aitken"
$ws.Cells.Item(19, 4).Value = 12960480531.96818
$ws.Cells.Item(19, 5).Value = 18
$ws.Cells.Item(19, 6).Value = 1.092768430709839
$ws.Cells.Item(19, 7).Value = 10974

$ws.Cells.Item(20, 2).Value = "bcsstk18.mtx"
$ws.Cells.Item(20, 3).Value = "MP"
$ws.Cells.Item(20, 4).Value = 42951395539.5179
$ws.Cells.Item(20, 5).Value = 39
$ws.Cells.Item(20, 6).Value = 2.845968008041382
$ws.Cells.Item(20, 7).Value = 11948

$ws.Cells.Item(21, 2).Value = "bcsstk18.mtx"
$ws.Cells.Item(21, 3).Value = "MP_Aitken"
$ws.Cells.Item(21, 4).Value = 42952170136.46172
$ws.Cells.Item(21, 5).Value = 32
$ws.Cells.Item(21, 6).Value = 2.339402914047241
$ws.Cells.Item(21, 7).Value = 11948

